$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "批准时间" (K1) is being replaced by a new column header "原编号".
# The other header cells downstream (L1 "保护现状", M1 "照片") keep their
# text; only the shared-string bookkeeping shifts as a side effect.
$ws.Range("K1").Value = "原编号"

# Move/save the active selection to K1, matching the edited workbook.
[void]$ws.Range("K1").Select()
